$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 0
    3 = 1
    4 = 1
    5 = 1
    6 = 0
    7 = 0
    8 = 0
    9 = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 1
    14 = 2
    15 = 1
    16 = 2
    17 = 1
    18 = 0
    19 = 3
    20 = 0
    21 = 0
    22 = 1
    23 = 1
    24 = 2
    25 = 1
    26 = 1
    27 = 0
    28 = 3
    29 = 0
    30 = 1
    31 = 1
    32 = 1
    33 = 0
    34 = 1
    35 = 1
    36 = 2
    37 = 1
    38 = 2
    39 = 0
    40 = 1
    41 = 2
    42 = 0
    43 = 0
    44 = 0
    45 = 1
    46 = 2
    47 = 1
    48 = 1
    49 = 1
    50 = 2
    51 = 1
    52 = 2
    53 = 1
    54 = 0
    55 = 1
    56 = 1
    57 = 2
    58 = 2
    59 = 2
    60 = 0
    61 = 2
    62 = 0
    63 = 2
    64 = 3
    65 = 3
    66 = 2
    67 = 1
    68 = 1
    69 = 1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}
